$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date bumped by 1 day (46078 -> 46079) for all data rows 2-9
$ws.Range("C2").Value = 46079
$ws.Range("C3").Value = 46079
$ws.Range("C4").Value = 46079
$ws.Range("C5").Value = 46079
$ws.Range("C6").Value = 46079
$ws.Range("C7").Value = 46079
$ws.Range("C8").Value = 46079
$ws.Range("C9").Value = 46079

# Rows 4-9 get reordered (Beteckning / Datum / Area columns) per the refreshed source data
$ws.Range("A4").Value = "A 17908-2021"
$ws.Range("B4").Value = 44301
$ws.Range("G4").Value = 0.9

$ws.Range("A5").Value = "A 54782-2022"
$ws.Range("B5").Value = 44883
$ws.Range("G5").Value = 5.5

$ws.Range("A6").Value = "A 25617-2024"
$ws.Range("B6").Value = 45463
$ws.Range("G6").Value = 2.3

$ws.Range("A7").Value = "A 45983-2023"
$ws.Range("B7").Value = 45196
$ws.Range("G7").Value = 0.6

$ws.Range("A8").Value = "A 843-2024"
$ws.Range("B8").Value = 45300
$ws.Range("G8").Value = 0.8

$ws.Range("A9").Value = "A 844-2024"
$ws.Range("B9").Value = 45300
$ws.Range("G9").Value = 1.2
